$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)
}

Replace-Text "English" "Английский"
Replace-Text " / Portuguese / French / Thai / Vietnamese / Spanish" " / португальский / французский / тайский / вьетнамский / испанский"
Replace-Text "Brief" "Кратко"
Replace-Text "An email sent to confirmed attendees of the event. We want to share the flight and accommodation booking details with them." "Письмо, отправленное подтверждённым участникам мероприятия. Мы хотим поделиться с ними деталями вашей брони на перелёт и проживание."
Replace-Text "Target audience" "Целевая аудитория"
Replace-Text "Event attendees" "Участники мероприятия"
Replace-Text "Subject: " "Тема: "
Replace-Text "Here are your booking details for" "Вот ваши детали бронирования для"
Replace-Text "We can’t wait to meet you! " "Очень ждем встречи с вами! "
Replace-Text "Hi " "Здравствуйте, "
Replace-Text "We hope you’re as excited as we are for " "Мы надеемся, что вы так же радостны, как и мы, по поводу "
Replace-Text ". As we’re nearing the event, we’ve made all the preparations to have you with us for this " ". Поскольку мы приближаемся к мероприятию, мы сделали все приготовления для того, чтобы вы смогли посетить этот "
Replace-Text "conference/seminar/trip" "конференцию/семинар/поездку"
Replace-Text "In this email, we’ve linked/attached the following documents:" "К этому письму приложены/прикреплены следующие документы:"
Replace-Text "Your return flight tickets" "Ваши билеты на обратный рейс"
Replace-Text "Your accommodation booking details" "Данные бронирования жилья"
Replace-Text "Your visa information " "Информация о вашей визе "
Replace-Text "(if applicable)" "(если применимо)"
Replace-Text "If you have any questions, please contact us via " "Если у вас есть вопросы, свяжитесь с нами через "
Replace-Text "live chat" "чат"
Replace-Text " or " " или "
Replace-Text "If you have any questions, please contact your country manager, " "Если у вас есть вопросы, пожалуйста, свяжитесь с вашим региональным менеджером, "
Replace-Text ", at " ", по адресу "
Replace-Text "See you on the " "Увидимся "
Replace-Text "[DD]th" "[DD]-го"

foreach ($c in $d.Comments) {
    $c.Range.Find.Execute("choose either one", $true, $true, $false, $false, $false, $true, 1, $false, "выберите один из", 2)
    $c.Range.Find.Execute("check if these are the documents included", $true, $true, $false, $false, $false, $true, 1, $false, "проверьте, включены ли эти документы", 2)
    $c.Range.Find.Execute("choose one", $true, $true, $false, $false, $false, $true, 1, $false, "выберите один", 2)
}
